# Weekly refresh: a new price observation is prepended to the historical
# series (dated 2023-03-29 / serial 45014), pushing the existing rows
# 66..186 down to 67..187.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 66; rows 66..186 shift down to 67..187.
$ws.Rows("66").Insert()

# Fill in the new weekly observation for Camote (Hortaliza) at Vega Modelo de Temuco.
$ws.Cells.Item(66, 1).Value = 10
$ws.Cells.Item(66, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(66, 3).Value = "La Araucanía"
$ws.Cells.Item(66, 4).Value = 45014
$ws.Cells.Item(66, 5).Value = 9
$ws.Cells.Item(66, 6).Value = 100114002
$ws.Cells.Item(66, 7).Value = "Camote"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 40
$ws.Cells.Item(66, 11).Value = 26000
$ws.Cells.Item(66, 12).Value = 26000
$ws.Cells.Item(66, 13).Value = 26000
$ws.Cells.Item(66, 14).Value = "$/malla 20 kilos"
$ws.Cells.Item(66, 15).Value = "Perú"
$ws.Cells.Item(66, 16).Value = 1300
$ws.Cells.Item(66, 17).Value = 20
$ws.Cells.Item(66, 18).Value = "Hortaliza"
